# Add "NETWORK VIRTUALIZATION CPU OVERHEAD" section (columns L:U) to the
# "summary" worksheet, mirroring the existing CASE A/B/D/"D in GUEST" CPU
# tables (columns A:J) but with the new total-virtualization-overhead figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New section title row (row 18) -----------------------------------
$ws.Range("A18").Value = "NETWORK VIRTUALIZATION CPU OVERHEAD"
$ws.Range("L18").Value = "TOTAL VIRTAUALIZATION CPU OVERHEAD"

# --- Column headers for each of the four sub-tables --------------------
# Each sub-table header row already has A:J filled in with the case name
# + GUEST/USER/NICE/SYS/IRQ/SOFTIRQ/STEAL/WAIT/IDLE labels (styled).
# We replicate the same row into L:U.
$headerRows = 19, 24, 29, 34

foreach ($r in $headerRows) {
    $caseLabel = $ws.Range("A$r").Value()
    $ws.Range("L$r").Value = $caseLabel

    $ws.Range("M$r").Value = "GUEST"
    $ws.Range("N$r").Value = "USER"
    $ws.Range("O$r").Value = "NICE"
    $ws.Range("P$r").Value = "SYS"
    $ws.Range("Q$r").Value = "IRQ"
    $ws.Range("R$r").Value = "SOFTIRQ"
    $ws.Range("S$r").Value = "STEAL"
    $ws.Range("T$r").Value = "WAIT"
    $ws.Range("U$r").Value = "IDLE"

    # Match the wrap-text / vertically-centered style used by B$r:J$r
    $ws.Range("B$r").Copy() | Out-Null
    $ws.Range("M$r`:U$r").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

# --- Data rows: GIGABIT / INFINIBAND for each sub-table -----------------
# Values: Guest, User, Nice, Sys, Irq, Softirq, Steal, Wait, Idle

$ws.Range("L20").Value = "GIGABIT"
$ws.Range("M20").Value = 5.7
$ws.Range("N20").Value = 5.7
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 68.25
$ws.Range("Q20").Value = 0
$ws.Range("R20").Value = 3.69
$ws.Range("S20").Value = 0
$ws.Range("T20").Value = 0
$ws.Range("U20").Value = 222.36

$ws.Range("L21").Value = "INFINIBAND"
$ws.Range("M21").Value = 68.61
$ws.Range("N21").Value = 68.64
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 50.7
$ws.Range("Q21").Value = 0
$ws.Range("R21").Value = 33.659999999999997
$ws.Range("S21").Value = 0
$ws.Range("T21").Value = 0
$ws.Range("U21").Value = 146.91

$ws.Range("L25").Value = "GIGABIT"
$ws.Range("M25").Value = 194.04
$ws.Range("N25").Value = 194.07
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = 3.09
$ws.Range("Q25").Value = 0
$ws.Range("R25").Value = 2.91
$ws.Range("S25").Value = 0
$ws.Range("T25").Value = 0
$ws.Range("U25").Value = 99.93

$ws.Range("L26").Value = "INFINIBAND"
$ws.Range("M26").Value = 125.28
$ws.Range("N26").Value = 125.31
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = 43.65
$ws.Range("Q26").Value = 0
$ws.Range("R26").Value = 31.08
$ws.Range("S26").Value = 0
$ws.Range("T26").Value = 0
$ws.Range("U26").Value = 99.96

$ws.Range("L30").Value = "GIGABIT"
$ws.Range("M30").Value = 199.83
$ws.Range("N30").Value = 199.89
$ws.Range("O30").Value = 0
$ws.Range("P30").Value = 0.24
$ws.Range("Q30").Value = 0
$ws.Range("R30").Value = 0
$ws.Range("S30").Value = 0
$ws.Range("T30").Value = 0
$ws.Range("U30").Value = 99.87

$ws.Range("L31").Value = "INFINIBAND"
$ws.Range("M31").Value = 199.74
$ws.Range("N31").Value = 199.77
$ws.Range("O31").Value = 0
$ws.Range("P31").Value = 0.24
$ws.Range("Q31").Value = 0
$ws.Range("R31").Value = 0
$ws.Range("S31").Value = 0
$ws.Range("T31").Value = 0
$ws.Range("U31").Value = 99.99

$ws.Range("L35").Value = "GIGABIT"
$ws.Range("M35").Value = 0
$ws.Range("N35").Value = 0.02
$ws.Range("O35").Value = 0
$ws.Range("P35").Value = 0.86
$ws.Range("Q35").Value = 0
$ws.Range("R35").Value = 2.78
$ws.Range("S35").Value = 0
$ws.Range("T35").Value = 0
$ws.Range("U35").Value = 196.34

$ws.Range("L36").Value = "INFINIBAND"
$ws.Range("M36").Value = 0
$ws.Range("N36").Value = 0.66
$ws.Range("O36").Value = 0
$ws.Range("P36").Value = 40.299999999999997
$ws.Range("Q36").Value = 0
$ws.Range("R36").Value = 41.42
$ws.Range("S36").Value = 0
$ws.Range("T36").Value = 0
$ws.Range("U36").Value = 117.64
